$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RC Band Pass Filter")

# Update component values: R2, C1, C2
$ws.Range("D2").Value = 2000
$ws.Range("D3").Value = 0.0000033
$ws.Range("D4").Value = 0.0000033

# Update the active selection on this sheet to D3
$ws.Range("D3").Select()

# Refresh chart so cached numCache values reflect the new data
$co = $ws.ChartObjects().Item(1)
$co.Chart.Refresh()
